$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3: E3 changes from 30 to 35 (D3 formula recalculates automatically)
$ws.Range("E3").Value = 35

# Row 4 (A4 date serial 42144 = May 21, 2019 under the workbook's 1904 date system)
$ws.Range("A4").Value = 42144
$ws.Range("B4").Value = 17
$ws.Range("C4").Value = 8
$ws.Range("D4").Value = 10
$ws.Range("E4").Value = 35

# Row 5
$ws.Range("A5").Value = 42151
$ws.Range("B5").Value = 20
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 10
$ws.Range("E5").Value = 35

# Row 6
$ws.Range("A6").Value = 42158
$ws.Range("B6").Value = 24
$ws.Range("C6").Value = 9
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 35

# Row 7
$ws.Range("A7").Value = 42165
$ws.Range("B7").Value = 28
$ws.Range("C7").Value = 6
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 35

# Row 8
$ws.Range("A8").Value = 42172
$ws.Range("B8").Value = 30
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 35

# Row 9 - "Final Release" label moves here from row 4
$ws.Range("A9").Value = "Final Release"
$ws.Range("B9").Value = 31
$ws.Range("C9").Value = 4
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 35

# Update selection to C9
$ws.Range("C9").Select()
